$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timing values (rows 5-7) ---
$ws.Range("B5").Value2 = 0.0003271102905273438
$ws.Range("B6").Value2 = 0.0003018379211425781
$ws.Range("B7").Value2 = 0.0004138946533203125

# --- Update embedding list text representation (row 8) ---
$ws.Range("A8").Value2 = "[[2, 1], [1, 2], [1, 1], [1, 0], [0, 1]]"

# --- Insert a new row before the old row 37 ("Movement times") ---
# This shifts rows 37-41 down to 38-42 and extends the used range to N42.
$ws.Rows.Item(37).EntireRow.Insert()

$ws.Range("A37").Value2 = "move_fidelity"
$ws.Range("B37").Value2 = 1

# --- Update the "total time:" value, now on row 41 after the insert ---
$ws.Range("B41").Value2 = 0.004017114639282227
